$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8040655475735718
$ws.Range("C2").Value = 0.1436306051655514
$ws.Range("E2").Value = 0.161725639684466
$ws.Range("F2").Value = 3.103742072394027
$ws.Range("G2").Value = 2.069768339182161
$ws.Range("H2").Value = 1.749238341856895
$ws.Range("J2").Value = 0.1520405910991647
$ws.Range("K2").Value = 0.3686519332727869
$ws.Range("L2").Value = 0.2948185936199366
$ws.Range("M2").Value = 0.2417527352607607
$ws.Range("B3").Value = 0.7769275376238625
$ws.Range("C3").Value = 0.142555060226357
$ws.Range("E3").Value = 0.16185835895884
$ws.Range("F3").Value = 3.095409000558973
$ws.Range("G3").Value = 2.064397557424357
$ws.Range("H3").Value = 1.751565992879023
$ws.Range("J3").Value = 0.1521971606541257
$ws.Range("K3").Value = 0.3434197027700066
$ws.Range("L3").Value = 0.29258987759588
$ws.Range("M3").Value = 0.2366164612222086
$ws.Range("B4").Value = 0.7607093242720566
$ws.Range("C4").Value = 0.1418764504728571
$ws.Range("E4").Value = 0.1619866915041595
$ws.Range("F4").Value = 3.091504985266866
$ws.Range("G4").Value = 2.061929026847992
$ws.Range("H4").Value = 1.753554716150646
$ws.Range("J4").Value = 0.1522930554011852
$ws.Range("K4").Value = 0.3281232047717566
$ws.Range("L4").Value = 0.291347270731805
$ws.Range("M4").Value = 0.2335879005035757
$ws.Range("B5").Value = 0.7542124659932199
$ws.Range("C5").Value = 0.1415953141558859
$ws.Range("E5").Value = 0.1620508043045152
$ws.Range("F5").Value = 3.090219198924729
$ws.Range("G5").Value = 2.061131493310711
$ws.Range("H5").Value = 1.75450591249637
$ws.Range("J5").Value = 0.1523320761234164
$ws.Range("K5").Value = 0.3219392885989549
$ws.Range("L5").Value = 0.2908726533729009
$ws.Range("M5").Value = 0.2323853140740475
$ws.Range("B6").Value = 0.7531404554676442
$ws.Range("C6").Value = 0.1415483535633726
$ws.Range("E6").Value = 0.1620621649588863
$ws.Range("F6").Value = 3.090024131161471
$ws.Range("G6").Value = 2.061011648659274
$ws.Range("H6").Value = 1.75467236265753
$ws.Range("J6").Value = 0.1523385521238101
$ws.Range("K6").Value = 0.3209154514399586
$ws.Range("L6").Value = 0.2907957647857558
$ws.Range("M6").Value = 0.2321875359534609
$ws.Range("B7").Value = 0.7606212505762926
$ws.Range("C7").Value = 0.1418726776046242
$ws.Range("E7").Value = 0.161987508258246
$ws.Range("F7").Value = 3.09148640893423
$ws.Range("G7").Value = 2.061917427294503
$ws.Range("H7").Value = 1.753566974229031
$ws.Range("J7").Value = 0.1522935818760063
$ws.Range("K7").Value = 0.3280396055035197
$ws.Range("L7").Value = 0.2913407411500231
$ws.Range("M7").Value = 0.2335715540015784
$ws.Range("B8").Value = 0.7946162991995607
$ws.Range("C8").Value = 0.1432635294076015
$ws.Range("E8").Value = 0.1617617012977259
$ws.Range("F8").Value = 3.100617300186997
$ws.Range("G8").Value = 2.067744346976113
$ws.Range("H8").Value = 1.749924827116075
$ws.Range("J8").Value = 0.1520946278697775
$ws.Range("K8").Value = 0.359911248616541
$ws.Range("L8").Value = 0.2940240818330935
$ws.Range("M8").Value = 0.2399558409478821
$ws.Range("B9").Value = 0.8647959528178717
$ws.Range("C9").Value = 0.1458471951402167
$ws.Range("E9").Value = 0.1616890231923307
$ws.Range("F9").Value = 3.12813756910252
$ws.Range("G9").Value = 2.085756302588749
$ws.Range("H9").Value = 1.747220038480918
$ws.Range("J9").Value = 0.151702428219485
$ws.Range("K9").Value = 0.4239632380815124
$ws.Range("L9").Value = 0.3002801582471335
$ws.Range("M9").Value = 0.2534643545996857
$ws.Range("B10").Value = 0.9184899858550182
$ws.Range("C10").Value = 0.1476589574974412
$ws.Range("E10").Value = 0.1618593324514457
$ws.Range("F10").Value = 3.1542156027011
$ws.Range("G10").Value = 2.103017058546271
$ws.Range("H10").Value = 1.74793667661649
$ws.Range("J10").Value = 0.1514128039121436
$ws.Range("K10").Value = 0.4719671459805568
$ws.Range("L10").Value = 0.3054774770083526
$ws.Range("M10").Value = 0.2639880926465068
$ws.Range("B11").Value = 0.9433783024253444
$ws.Range("C11").Value = 0.148464667519427
$ws.Range("E11").Value = 0.1619849880583111
$ws.Range("F11").Value = 3.167351274167672
$ws.Range("G11").Value = 2.111746919377822
$ws.Range("H11").Value = 1.74884965506331
$ws.Range("J11").Value = 0.1512806788369128
$ws.Range("K11").Value = 0.4940106969226008
$ws.Range("L11").Value = 0.3079713690786008
$ws.Range("M11").Value = 0.2689049726281638
$ws.Range("B12").Value = 0.9528691080327008
$ws.Range("C12").Value = 0.1487671316510699
$ws.Range("E12").Value = 0.1620394638286555
$ws.Range("F12").Value = 3.172508290650725
$ws.Range("G12").Value = 2.115179092754687
$ws.Range("H12").Value = 1.74927974345718
$ws.Range("J12").Value = 0.1512305897660244
$ws.Range("K12").Value = 0.5023875811418748
$ws.Range("L12").Value = 0.3089342787951352
$ws.Range("M12").Value = 0.2707854121763518
$ws.Range("B13").Value = 0.9508221577301015
$ws.Range("C13").Value = 0.1487021077956356
$ws.Range("E13").Value = 0.1620274255071728
$ws.Range("F13").Value = 3.171389506677457
$ws.Range("G13").Value = 2.114434290645363
$ws.Range("H13").Value = 1.749183364826706
$ws.Range("J13").Value = 0.1512413798799139
$ws.Range("K13").Value = 0.5005821617278912
$ws.Range("L13").Value = 0.3087260767175763
$ws.Range("M13").Value = 0.270379603757199
$ws.Range("B14").Value = 0.9441577933613416
$ws.Range("C14").Value = 0.1484896042634105
$ws.Range("E14").Value = 0.1619893318759011
$ws.Range("F14").Value = 3.167771882016893
$ws.Range("G14").Value = 2.112026753071405
$ws.Range("H14").Value = 1.74888334817399
$ws.Range("J14").Value = 0.1512765591187302
$ws.Range("K14").Value = 0.4946992787477598
$ws.Range("L14").Value = 0.3080502175622399
$ws.Range("M14").Value = 0.2690593069008287
$ws.Range("B15").Value = 0.9400842797858502
$ws.Range("C15").Value = 0.1483590962601085
$ws.Range("E15").Value = 0.161966895007339
$ws.Range("F15").Value = 3.165579785988029
$ws.Range("G15").Value = 2.110568527032541
$ws.Range("H15").Value = 1.748710564628198
$ws.Range("J15").Value = 0.1512981000315055
$ws.Range("K15").Value = 0.4910996749999299
$ws.Range("L15").Value = 0.3076386438062571
$ws.Range("M15").Value = 0.2682529960656765
$ws.Range("B16").Value = 0.9168727521271478
$ws.Range("C16").Value = 0.1476059320878704
$ws.Range("E16").Value = 0.1618520870355233
$ws.Range("F16").Value = 3.153382748968653
$ws.Range("G16").Value = 2.102464217433067
$ws.Range("H16").Value = 1.747888819082135
$ws.Range("J16").Value = 0.1514214308673374
$ws.Range("K16").Value = 0.4705306828280413
$ws.Range("L16").Value = 0.3053170938566581
$ws.Range("M16").Value = 0.2636693610567633
$ws.Range("B17").Value = 0.9027514473518465
$ws.Range("C17").Value = 0.1471391704717249
$ws.Range("E17").Value = 0.1617939709560972
$ws.Range("F17").Value = 3.146226086864473
$ws.Range("G17").Value = 2.097717418586342
$ws.Range("H17").Value = 1.747535012563304
$ws.Range("J17").Value = 0.151496992883124
$ws.Range("K17").Value = 0.4579649756209108
$ws.Range("L17").Value = 0.3039260146399982
$ws.Range("M17").Value = 0.2608905595307291
$ws.Range("B18").Value = 0.8946728142539087
$ws.Range("C18").Value = 0.1468689621626709
$ws.Range("E18").Value = 0.1617650812752061
$ws.Range("F18").Value = 3.142229562136208
$ws.Range("G18").Value = 2.095069805673347
$ws.Range("H18").Value = 1.747386763905837
$ws.Range("J18").Value = 0.1515404191358591
$ws.Range("K18").Value = 0.4507569460018601
$ws.Range("L18").Value = 0.3031381083176967
$ws.Range("M18").Value = 0.259304469448395
$ws.Range("B19").Value = 0.891945022533946
$ws.Range("C19").Value = 0.1467771750729128
$ws.Range("E19").Value = 0.1617560800753353
$ws.Range("F19").Value = 3.140896988511372
$ws.Range("G19").Value = 2.094187555583375
$ws.Range("H19").Value = 1.747346060289544
$ws.Range("J19").Value = 0.1515551165809397
$ws.Range("K19").Value = 0.4483197747775591
$ws.Range("L19").Value = 0.3028734364765313
$ws.Range("M19").Value = 0.2587695458348946
$ws.Range("B20").Value = 0.9042501780407122
$ws.Range("C20").Value = 0.1471890379511009
$ws.Range("E20").Value = 0.1617996881931454
$ws.Range("F20").Value = 3.146975528427959
$ws.Range("G20").Value = 2.098214171594947
$ws.Range("H20").Value = 1.74756695785706
$ws.Range("J20").Value = 0.151488952819701
$ws.Range("K20").Value = 0.4593006068419356
$ws.Range("L20").Value = 0.3040728349595554
$ws.Range("M20").Value = 0.2611851055219816
$ws.Range("B21").Value = 0.9461134867514716
$ws.Range("C21").Value = 0.1485520932502311
$ws.Range("E21").Value = 0.1620003341156462
$ws.Range("F21").Value = 3.168829505761366
$ws.Range("G21").Value = 2.11273047468552
$ws.Range("H21").Value = 1.748969181083311
$ws.Range("J21").Value = 0.1512662276667394
$ws.Range("K21").Value = 0.4964264258634898
$ws.Range("L21").Value = 0.3082482318394142
$ws.Range("M21").Value = 0.2694466082294866
$ws.Range("B22").Value = 0.9738589186052877
$ws.Range("C22").Value = 0.149427543370912
$ws.Range("E22").Value = 0.162171627456658
$ws.Range("F22").Value = 3.184177890247426
$ws.Range("G22").Value = 2.122954327241217
$ws.Range("H22").Value = 1.750377329151405
$ws.Range("J22").Value = 0.1511203354336894
$ws.Range("K22").Value = 0.5208619765521973
$ws.Range("L22").Value = 0.3110850334354751
$ws.Range("M22").Value = 0.2749538985496329
$ws.Range("B23").Value = 0.9590155273807284
$ws.Range("C23").Value = 0.1489617018883749
$ws.Range("E23").Value = 0.1620765420991468
$ws.Range("F23").Value = 3.175888730954114
$ws.Range("G23").Value = 2.117430222092935
$ws.Range("H23").Value = 1.749580796260773
$ws.Range("J23").Value = 0.1511982317013354
$ws.Range("K23").Value = 0.5078046199649862
$ws.Range("L23").Value = 0.3095611389099133
$ws.Range("M23").Value = 0.2720047152596266
$ws.Range("B24").Value = 0.9035724776048824
$ws.Range("C24").Value = 0.147166498662969
$ws.Range("E24").Value = 0.1617970893424712
$ws.Range("F24").Value = 3.146636338512309
$ws.Range("G24").Value = 2.097989336051512
$ws.Range("H24").Value = 1.747552343557601
$ws.Range("J24").Value = 0.1514925877788569
$ws.Range("K24").Value = 0.458696717608035
$ws.Range("L24").Value = 0.3040064205941775
$ws.Range("M24").Value = 0.2610519055127511
$ws.Range("B25").Value = 0.845435197793222
$ws.Range("C25").Value = 0.1451635149889157
$ws.Range("E25").Value = 0.1616692604467254
$ws.Range("F25").Value = 3.11966354989039
$ws.Range("G25").Value = 2.080177262480404
$ws.Range("H25").Value = 1.747476837249948
$ws.Range("J25").Value = 0.1518087719235366
$ws.Range("K25").Value = 0.4064694716847725
$ws.Range("L25").Value = 0.2984818253644832
$ws.Range("M25").Value = 0.2497044719657779

Write-Host "applied 380 kV results"
